# Auto-generated PowerShell Excel COM-interop script
# Updates loading_percent values for rows 2-25, columns C,D,E,F,G,H,J,L,O
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.62572784085251
$ws.Range("D2").Value = 8.862755317428553
$ws.Range("E2").Value = 14.33113577209335
$ws.Range("F2").Value = 37.46470187668108
$ws.Range("G2").Value = 42.83247596206232
$ws.Range("H2").Value = 17.26068741337964
$ws.Range("J2").Value = 10.70174230392222
$ws.Range("L2").Value = 9.757431744702671
$ws.Range("O2").Value = 28.31248297128687

$ws.Range("C3").Value = 13.59124863422232
$ws.Range("D3").Value = 8.873921610170257
$ws.Range("E3").Value = 14.34229924205262
$ws.Range("F3").Value = 37.26013579872049
$ws.Range("G3").Value = 42.33423386470002
$ws.Range("H3").Value = 17.23884865881928
$ws.Range("J3").Value = 10.71951035084092
$ws.Range("L3").Value = 9.775400896148057
$ws.Range("O3").Value = 28.18243555185088

$ws.Range("C4").Value = 13.57300577819367
$ws.Range("D4").Value = 8.881897797227888
$ws.Range("E4").Value = 14.35144590203804
$ws.Range("F4").Value = 37.14330426916216
$ws.Range("G4").Value = 42.03776962370883
$ws.Range("H4").Value = 17.22888948356966
$ws.Range("J4").Value = 10.73176849844094
$ws.Range("L4").Value = 9.787315083032475
$ws.Range("O4").Value = 28.10908694031794

$ws.Range("C5").Value = 13.56631267829636
$ws.Range("D5").Value = 8.885430033876764
$ws.Range("E5").Value = 14.35574987205515
$ws.Range("F5").Value = 37.09793849220998
$ws.Range("G5").Value = 41.9194764728114
$ws.Range("H5").Value = 17.22570137271192
$ws.Range("J5").Value = 10.73710307738097
$ws.Range("L5").Value = 9.792392179752357
$ws.Range("O5").Value = 28.08085364041314

$ws.Range("C6").Value = 13.56524617732769
$ws.Range("D6").Value = 8.886033589904891
$ws.Range("E6").Value = 14.35649937399978
$ws.Range("F6").Value = 37.09054211389156
$ws.Range("G6").Value = 41.89999007353505
$ws.Range("H6").Value = 17.22522462595816
$ws.Range("J6").Value = 10.73800937736525
$ws.Range("L6").Value = 9.793248646511318
$ws.Range("O6").Value = 28.07626619277322

$ws.Range("C7").Value = 13.57291250648725
$ws.Range("D7").Value = 8.88194429261339
$ws.Range("E7").Value = 14.35150161192673
$ws.Range("F7").Value = 37.14268331700421
$ws.Range("G7").Value = 42.03616390314618
$ws.Range("H7").Value = 17.22884296042276
$ws.Range("J7").Value = 10.73183906845057
$ws.Range("L7").Value = 9.787382655234669
$ws.Range("O7").Value = 28.10869944022895

$ws.Range("C8").Value = 13.61323493729683
$ws.Range("D8").Value = 8.866373138898711
$ws.Range("E8").Value = 14.3345093573034
$ws.Range("F8").Value = 37.39237054819648
$ws.Range("G8").Value = 42.65880465962623
$ws.Range("H8").Value = 17.25244325448842
$ws.Range("J8").Value = 10.70758896840129
$ws.Range("L8").Value = 9.763444897773969
$ws.Range("O8").Value = 28.26630730001701

$ws.Range("C9").Value = 13.71529058171047
$ws.Range("D9").Value = 8.844714389750035
$ws.Range("E9").Value = 14.31936133740965
$ws.Range("F9").Value = 37.94978392962692
$ws.Range("G9").Value = 43.94845170283354
$ws.Range("H9").Value = 17.32594814660332
$ws.Range("J9").Value = 10.67072554505185
$ws.Range("L9").Value = 9.723475175503429
$ws.Range("O9").Value = 28.62588138305022

$ws.Range("C10").Value = 13.80390930595297
$ws.Range("D10").Value = 8.834197609021
$ws.Range("E10").Value = 14.31928596505051
$ws.Range("F10").Value = 38.39796505636173
$ws.Range("G10").Value = 44.92898426115374
$ws.Range("H10").Value = 17.39631425799097
$ws.Range("J10").Value = 10.65014734746838
$ws.Range("L10").Value = 9.698334533188333
$ws.Range("O10").Value = 28.91926276991265

$ws.Range("C11").Value = 13.84709189678689
$ws.Range("D11").Value = 8.83058104224229
$ws.Range("E11").Value = 14.32164344846577
$ws.Range("F11").Value = 38.6096184240628
$ws.Range("G11").Value = 45.38031206323657
$ws.Range("H11").Value = 17.43181581195756
$ws.Range("J11").Value = 10.64219537716486
$ws.Range("L11").Value = 9.687809496925768
$ws.Range("O11").Value = 29.05869084702694

$ws.Range("C12").Value = 13.86384765265519
$ws.Range("D12").Value = 8.829379038811565
$ws.Range("E12").Value = 14.32287906416177
$ws.Range("F12").Value = 38.69082727827905
$ws.Range("G12").Value = 45.55181696055553
$ws.Range("H12").Value = 17.44575508117034
$ws.Range("J12").Value = 10.63938650960057
$ws.Range("L12").Value = 9.683954600182055
$ws.Range("O12").Value = 29.11231200333798

$ws.Range("C13").Value = 13.86022121028751
$ws.Range("D13").Value = 8.829630468129297
$ws.Range("E13").Value = 14.32259772063713
$ws.Range("F13").Value = 38.67329127754634
$ws.Range("G13").Value = 45.51485633245254
$ws.Range("H13").Value = 17.44273107786235
$ws.Range("J13").Value = 10.63998245347786
$ws.Range("L13").Value = 9.684779013963956
$ws.Range("O13").Value = 29.10072771395284

$ws.Range("C14").Value = 13.84846237205786
$ws.Range("D14").Value = 8.830478797838284
$ws.Range("E14").Value = 14.32173823791998
$ws.Range("F14").Value = 38.61627861159096
$ws.Range("G14").Value = 45.39441083596205
$ws.Range("H14").Value = 17.43295271075668
$ws.Range("J14").Value = 10.64196023573304
$ws.Range("L14").Value = 9.687489734858824
$ws.Range("O14").Value = 29.06308601017928

$ws.Range("C15").Value = 13.84131200150168
$ws.Range("D15").Value = 8.83102022725474
$ws.Range("E15").Value = 14.32125640013919
$ws.Range("F15").Value = 38.58149305073956
$ws.Range("G15").Value = 45.32070743187458
$ws.Range("H15").Value = 17.4270275175376
$ws.Range("J15").Value = 10.64319803084871
$ws.Range("L15").Value = 9.689167140311811
$ws.Range("O15").Value = 29.04013545344816

$ws.Range("C16").Value = 13.80114419290252
$ws.Range("D16").Value = 8.834457422766041
$ws.Range("E16").Value = 14.31917994096706
$ws.Range("F16").Value = 38.38428456588805
$ws.Range("G16").Value = 44.89958124957404
$ws.Range("H16").Value = 17.39406386747872
$ws.Range("J16").Value = 10.65069536131443
$ws.Range("L16").Value = 9.699040679950883
$ws.Range("O16").Value = 28.91026791198709

$ws.Range("C17").Value = 13.77723115553592
$ws.Range("D17").Value = 8.836864812886517
$ws.Range("E17").Value = 14.31851799987657
$ws.Range("F17").Value = 38.26525429078877
$ws.Range("G17").Value = 44.64247014623285
$ws.Range("H17").Value = 17.37473160250477
$ws.Range("J17").Value = 10.6556554701805
$ws.Range("L17").Value = 9.705330983097641
$ws.Range("O17").Value = 28.83210214360567

$ws.Range("C18").Value = 13.763747689522
$ws.Range("D18").Value = 8.838359402522299
$ws.Range("E18").Value = 14.31836246696413
$ws.Range("F18").Value = 38.19752649068062
$ws.Range("G18").Value = 44.4950931443814
$ws.Range("H18").Value = 17.36394127360017
$ws.Range("J18").Value = 10.65864104749174
$ws.Range("L18").Value = 9.709034822654331
$ws.Range("O18").Value = 28.78770716343045

$ws.Range("C19").Value = 13.75922917870525
$ws.Range("D19").Value = 8.838884335232633
$ws.Range("E19").Value = 14.31834851326142
$ws.Range("F19").Value = 38.17472302161068
$ws.Range("G19").Value = 44.44528578053129
$ws.Range("H19").Value = 17.36034456122506
$ws.Range("J19").Value = 10.65967470367349
$ws.Range("L19").Value = 9.710303630838956
$ws.Range("O19").Value = 28.77277368737875

$ws.Range("C20").Value = 13.77974878931563
$ws.Range("D20").Value = 8.836597168198072
$ws.Range("E20").Value = 14.31856516455171
$ws.Range("F20").Value = 38.27784957978065
$ws.Range("G20").Value = 44.66978882063506
$ws.Range("H20").Value = 17.37675554200241
$ws.Range("J20").Value = 10.65511373079685
$ws.Range("L20").Value = 9.704652489920079
$ws.Range("O20").Value = 28.8403649125191

$ws.Range("C21").Value = 13.85190535474515
$ws.Range("D21").Value = 8.830225079856739
$ws.Range("E21").Value = 14.32198139212292
$ws.Range("F21").Value = 38.6329963051614
$ws.Range("G21").Value = 45.42977368428577
$ws.Range("H21").Value = 17.43581145843749
$ws.Range("J21").Value = 10.6413738230912
$ws.Range("L21").Value = 9.686689986092652
$ws.Range("O21").Value = 29.07412025053699

$ws.Range("C22").Value = 13.90141074261738
$ws.Range("D22").Value = 8.827036734603782
$ws.Range("E22").Value = 14.32621210263082
$ws.Range("F22").Value = 38.87125753551988
$ws.Range("G22").Value = 45.92987274775989
$ws.Range("H22").Value = 17.47729322848912
$ws.Range("J22").Value = 10.63357343013889
$ws.Range("L22").Value = 9.67571213245432
$ws.Range("O22").Value = 29.23167117122135

$ws.Range("C23").Value = 13.87477723634493
$ws.Range("D23").Value = 8.828649230142746
$ws.Range("E23").Value = 14.32377166416428
$ws.Range("F23").Value = 38.74354959853644
$ws.Range("G23").Value = 45.66270203397099
$ws.Range("H23").Value = 17.4548919191058
$ws.Range("J23").Value = 10.63762882299486
$ws.Range("L23").Value = 9.681501649785226
$ws.Range("O23").Value = 29.14715812794682

$ws.Range("C24").Value = 13.77860974266244
$ws.Range("D24").Value = 8.836717825994858
$ws.Range("E24").Value = 14.31854314040871
$ws.Range("F24").Value = 38.27215305146698
$ws.Range("G24").Value = 44.65743666625356
$ws.Range("H24").Value = 17.37583950907082
$ws.Range("J24").Value = 10.65535823393429
$ws.Range("L24").Value = 9.704958964230487
$ws.Range("O24").Value = 28.83662761886628

$ws.Range("C25").Value = 13.6852569760637
$ws.Range("D25").Value = 8.84962483428802
$ws.Range("E25").Value = 14.32151563862596
$ws.Range("F25").Value = 37.79199649423816
$ws.Range("G25").Value = 43.59311989568452
$ws.Range("H25").Value = 17.30317030201023
$ws.Range("J25").Value = 10.67955469390842
$ws.Range("L25").Value = 9.733544246812029
$ws.Range("O25").Value = 28.52335573645837
